$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.706.62'
$ws.Range("E2").Value = '  -0.99%  '

$ws.Range("D3").Value = '1.598.16'
$ws.Range("E3").Value = '  -1.31%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '''211.44'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.82%  '

$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -1.18%  '

$ws.Range("E9").Value = '  -1.70%  '

$ws.Range("D10").Value = '''19.74'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.83%  '

$ws.Range("D11").Value = '''0.0837'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.01%  '

$ws.Range("D12").Value = '1.822.79'
$ws.Range("E12").Value = '  -1.30%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.03'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.69%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.560.82'
$ws.Range("E14").Value = '  -3.66%  '

$ws.Range("D15").Value = '''0.523'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.17%  '

$ws.Range("D16").Value = '''65.14'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.64%  '

$ws.Range("D17").Value = '26.707.00'
$ws.Range("E17").Value = '  -0.98%  '

$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("E18").Value = '  -0.65%  '

$ws.Range("D19").Value = '''209.95'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.45%  '

$ws.Range("E20").Value = '  -0.05%  '

$ws.Range("D21").Value = '''6.73'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.84%  '

$ws.Range("E22").Value = '  -1.13%  '

$ws.Range("D23").Value = '''2.30'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.58%  '

$ws.Range("E24").Value = '  -0.16%  '

$ws.Range("D25").Value = '''146.75'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("E27").Value = '  -3.89%  '

$ws.Range("E28").Value = '  +1.29%  '

$ws.Range("E29").Value = '  -0.97%  '

$ws.Range("D30").Value = '''0.0504'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.50%  '

$ws.Range("D31").Value = '''1.15'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.04%  '

$ws.Range("E32").Value = '  -1.85%  '

$ws.Range("D33").Value = '''0.664'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.50%  '

$ws.Range("E34").Value = '  -2.58%  '

$ws.Range("D35").Value = '1.298.82'
$ws.Range("E35").Value = '  -3.09%  '

$ws.Range("E37").Value = '  -4.81%  '

$ws.Range("E38").Value = '  -1.85%  '

$ws.Range("D39").Value = '''0.844'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.63%  '

$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("E41").Value = '  -0.54%  '

$ws.Range("E42").Value = '  +0.79%  '

$ws.Range("E43").Value = '  -0.81%  '

$ws.Range("D44").Value = '''63.74'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.09%  '

$ws.Range("D45").Value = '1.734.77'
$ws.Range("E45").Value = '  -1.35%  '

$ws.Range("D46").Value = '''90.14'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.48%  '

$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = '''0.877'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +9.84%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''1.63'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.37%  '

$ws.Range("D49").Value = '''0.0989'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.32%  '

$ws.Range("E50").Value = '  -2.05%  '

$ws.Range("D51").Value = '''7.50'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.63%  '

Write-Output "Applied crypto price/volume updates"